$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "CD4051 select A" config-param note, replacing the old IR-input /
# indicator-LED notes on GPIO3-GPIO6, plus a "*" marker on GPIO2 and the
# IR-input note moved down to GPIO10 (row 14).
$ws.Range("H5").Value = "*"
$ws.Range("E6").Value = "CD4051  select A"
$ws.Range("E7").Value = "CD4051  select A"
$ws.Range("E8").Value = "CD4051  select A"
$ws.Range("E9").ClearContents()
$ws.Range("E14").Value = "(IR input)"

# Move the selection / scroll position shown when the sheet is reopened.
$ws.Range("E10").Select()

# Page setup: landscape, smaller paper, tighter margins, gridlines on.
$ps = $ws.PageSetup
$ps.PrintGridlines = $true
$ps.Orientation = 2
$ps.PaperSize = 9
$ps.BlackAndWhite = $true
$ps.LeftMargin = 22.677165354330707
$ps.RightMargin = 22.677165354330707
$ps.TopMargin = 28.34645669291339
$ps.BottomMargin = 53.858267716535432
$ps.HeaderMargin = 28.34645669291339
$ps.FooterMargin = 22.677165354330707
